$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 1666.6666
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 1666.6666
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 4999.9998
$ws.Range("N48").Value = -5583.9998
$ws.Range("M48").ClearContents()
$ws.Range("H51").Value = 5688.9565
$ws.Range("I51").Value = 2333.3333
$ws.Range("J51").Value = 6192.3
$ws.Range("K51").Value = 2333.3333
$ws.Range("L51").Value = 6192.3
$ws.Range("M51").Value = -1849.3333
$ws.Range("N51").Value = -7160.3
$ws.Range("H56").Value = 1666.6666
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 1666.6666
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 4999.9998
$ws.Range("N56").Value = -6067.9998
$ws.Range("M56").ClearContents()
$ws.Range("H112").Value = 1299.4546
$ws.Range("J112").Value = 1279.4
$ws.Range("L112").Value = 3838.2
$ws.Range("N112").Value = -6054.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3908.1316
$ws.Range("I32").Value = 2467.853
$ws.Range("J32").Value = 16150.5
$ws.Range("K32").Value = 2467.853
$ws.Range("L32").Value = 16150.5
$ws.Range("M32").Value = -2180.853
$ws.Range("N32").Value = -16724.5
$ws.Range("H45").Value = 1271.3914
$ws.Range("I45").Value = 988.7692
$ws.Range("J45").Value = 1638.8
$ws.Range("K45").Value = 988.7692
$ws.Range("L45").Value = 1638.8
$ws.Range("M45").Value = -611.7692
$ws.Range("N45").Value = -2392.8
$ws.Range("H61").Value = 2442.25
$ws.Range("I61").Value = 586.625
$ws.Range("K61").Value = 586.625
$ws.Range("M61").Value = -374.625
$ws.Range("H88").Value = 13434
$ws.Range("I88").Value = 18484.334
$ws.Range("J88").Value = 3333.3333
$ws.Range("K88").Value = 18484.334
$ws.Range("L88").Value = 3333.3333
$ws.Range("M88").Value = -18078.334
$ws.Range("N88").Value = -4145.3333
$ws.Range("H91").Value = 13434
$ws.Range("I91").Value = 18484.334
$ws.Range("J91").Value = 3333.3333
$ws.Range("K91").Value = 18484.334
$ws.Range("L91").Value = 3333.3333
$ws.Range("M91").Value = -17080.334
$ws.Range("N91").Value = -6141.3333
$ws.Range("H110").Value = 40924.266
$ws.Range("I110").Value = 47002.92
$ws.Range("J110").Value = 1413
$ws.Range("K110").Value = 47002.92
$ws.Range("L110").Value = 1413
$ws.Range("M110").Value = -44957.92
$ws.Range("N110").Value = -5503
$ws.Range("H122").Value = 1807.4117
$ws.Range("I122").Value = 959.63635
$ws.Range("J122").Value = 3361.6667
$ws.Range("K122").Value = 2878.90905
$ws.Range("L122").Value = 10085.0001
$ws.Range("M122").Value = -428.9090500000002
$ws.Range("N122").Value = -14985.0001
$ws.Range("H136").Value = 2442.25
$ws.Range("I136").Value = 586.625
$ws.Range("K136").Value = 1759.875
$ws.Range("M136").Value = 790.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2252.6667
$ws.Range("I86").Value = 2569.1
$ws.Range("J86").Value = 1619.8
$ws.Range("K86").Value = 2569.1
$ws.Range("L86").Value = 1619.8
$ws.Range("M86").Value = -1446.1
$ws.Range("N86").Value = -3865.8
$ws.Range("H89").Value = 2252.6667
$ws.Range("I89").Value = 2569.1
$ws.Range("J89").Value = 1619.8
$ws.Range("K89").Value = 12845.5
$ws.Range("L89").Value = 8099
$ws.Range("M89").Value = -7229.5
$ws.Range("N89").Value = -19331
$ws.Range("H94").Value = 801.9583
$ws.Range("I94").Value = 737.85
$ws.Range("J94").Value = 1122.5
$ws.Range("K94").Value = 737.85
$ws.Range("L94").Value = 1122.5
$ws.Range("M94").Value = -286.85
$ws.Range("N94").Value = -2024.5
$ws.Range("H134").Value = 2582.8438
$ws.Range("I134").Value = 1232.6875
$ws.Range("J134").Value = 3933
$ws.Range("K134").Value = 3698.0625
$ws.Range("L134").Value = 11799
$ws.Range("M134").Value = -1163.0625
$ws.Range("N134").Value = -16869

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2100.639
$ws.Range("I58").Value = 1504.1428
$ws.Range("J58").Value = 2480.2273
$ws.Range("K58").Value = 1504.1428
$ws.Range("L58").Value = 2480.2273
$ws.Range("M58").Value = -1301.1428
$ws.Range("N58").Value = -2886.2273
$ws.Range("H122").Value = 715337.7
$ws.Range("I122").Value = 1602.8
$ws.Range("J122").Value = 1111857.1
$ws.Range("K122").Value = 4808.4
$ws.Range("L122").Value = 3335571.3
$ws.Range("M122").Value = -2358.4
$ws.Range("N122").Value = -3340471.3
$ws.Range("H136").Value = 2100.639
$ws.Range("I136").Value = 1504.1428
$ws.Range("J136").Value = 2480.2273
$ws.Range("K136").Value = 4512.428400000001
$ws.Range("L136").Value = 7440.6819
$ws.Range("M136").Value = -1962.428400000001
$ws.Range("N136").Value = -12540.6819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1362.76
$ws.Range("I5").Value = 754.73334
$ws.Range("J5").Value = 2274.8
$ws.Range("K5").Value = 2264.20002
$ws.Range("L5").Value = 6824.400000000001
$ws.Range("M5").Value = -2152.20002
$ws.Range("N5").Value = -7048.400000000001
$ws.Range("H87").Value = 8858.5
$ws.Range("I87").Value = 6004.6665
$ws.Range("J87").Value = 17420
$ws.Range("K87").Value = 18013.9995
$ws.Range("L87").Value = 52260
$ws.Range("M87").Value = -16765.9995
$ws.Range("N87").Value = -54756
$ws.Range("H90").Value = 8858.5
$ws.Range("I90").Value = 6004.6665
$ws.Range("J90").Value = 17420
$ws.Range("K90").Value = 54041.9985
$ws.Range("L90").Value = 156780
$ws.Range("M90").Value = -47801.9985
$ws.Range("N90").Value = -169260
$ws.Range("H132").Value = 973.7778
$ws.Range("I132").Value = 727.2727
$ws.Range("K132").Value = 6545.454299999999
$ws.Range("M132").Value = -4015.454299999999
$ws.Range("H135").Value = 1362.76
$ws.Range("I135").Value = 754.73334
$ws.Range("J135").Value = 2274.8
$ws.Range("K135").Value = 6792.60006
$ws.Range("L135").Value = 20473.2
$ws.Range("M135").Value = -4257.60006
$ws.Range("N135").Value = -25543.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 16333.333
$ws.Range("I39").Value = 16000
$ws.Range("J39").Value = 16500
$ws.Range("K39").Value = 16000
$ws.Range("L39").Value = 16500
$ws.Range("M39").Value = -15468
$ws.Range("N39").Value = -17564
$ws.Range("H41").Value = 1600
$ws.Range("J41").Value = 1900
$ws.Range("L41").Value = 1900
$ws.Range("N41").Value = -2610
$ws.Range("H132").Value = 5300.5835
$ws.Range("I132").Value = 2521
$ws.Range("J132").Value = 7286
$ws.Range("K132").Value = 7563
$ws.Range("L132").Value = 21858
$ws.Range("M132").Value = -5033
$ws.Range("N132").Value = -26918

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47352.863
$ws.Range("I7").Value = 78284.92
$ws.Range("J7").Value = 2673.2222
$ws.Range("K7").Value = 78284.92
$ws.Range("L7").Value = 2673.2222
$ws.Range("M7").Value = -78172.92
$ws.Range("N7").Value = -2897.2222
$ws.Range("H68").Value = 2709.875
$ws.Range("I68").Value = 2125
$ws.Range("J68").Value = 3294.75
$ws.Range("K68").Value = 2125
$ws.Range("L68").Value = 3294.75
$ws.Range("M68").Value = -1376
$ws.Range("N68").Value = -4792.75
$ws.Range("H71").Value = 2709.875
$ws.Range("I71").Value = 2125
$ws.Range("J71").Value = 3294.75
$ws.Range("K71").Value = 10625
$ws.Range("L71").Value = 16473.75
$ws.Range("M71").Value = -6881
$ws.Range("N71").Value = -23961.75
$ws.Range("H122").Value = 4277023
$ws.Range("I122").Value = 15874478
$ws.Range("J122").Value = 4276.4736
$ws.Range("K122").Value = 47623434
$ws.Range("L122").Value = 12829.4208
$ws.Range("M122").Value = -47620984
$ws.Range("N122").Value = -17729.4208
$ws.Range("H126").Value = 47352.863
$ws.Range("I126").Value = 78284.92
$ws.Range("J126").Value = 2673.2222
$ws.Range("K126").Value = 234854.76
$ws.Range("L126").Value = 8019.6666
$ws.Range("M126").Value = -232384.76
$ws.Range("N126").Value = -12959.6666
$ws.Range("H132").Value = 11214.556
$ws.Range("I132").Value = 16570.715
$ws.Range("J132").Value = 5446.385
$ws.Range("K132").Value = 49712.145
$ws.Range("L132").Value = 16339.155
$ws.Range("M132").Value = -47182.145
$ws.Range("N132").Value = -21399.155
$ws.Range("H136").Value = 22226028
$ws.Range("I136").Value = 4613.8
$ws.Range("J136").Value = 33336736
$ws.Range("K136").Value = 13841.4
$ws.Range("L136").Value = 100010208
$ws.Range("M136").Value = -11291.4
$ws.Range("N136").Value = -100015308

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 51000
$ws.Range("J46").Value = 51000
$ws.Range("L46").Value = 51000
$ws.Range("N46").Value = -51462
$ws.Range("H126").Value = 49032.332
$ws.Range("I126").Value = 72327.07000000001
$ws.Range("J126").Value = 2442.8572
$ws.Range("K126").Value = 216981.21
$ws.Range("L126").Value = 7328.571599999999
$ws.Range("M126").Value = -214511.21
$ws.Range("N126").Value = -12268.5716
$ws.Range("H134").Value = 51000
$ws.Range("J134").Value = 51000
$ws.Range("L134").Value = 153000
$ws.Range("N134").Value = -158070
$ws.Range("H136").Value = 7427.48
$ws.Range("I136").Value = 1309
$ws.Range("J136").Value = 10306.765
$ws.Range("K136").Value = 3927
$ws.Range("L136").Value = 30920.295
$ws.Range("M136").Value = -1377
$ws.Range("N136").Value = -36020.295

Write-Host "Applied all profit sheet updates"